# "Generate Report for Handback" -- refresh the handback-status report with
# the results of a new handoff/handback cycle:
#   - old source file 4032afde-8e5a-4dd2-a7f0-f64482d5e798.md -> c7e37e7b-0198-4224-a154-a931b3cbc493.md
#   - old source file 961d70ea-fabc-4d5e-9a78-48ffc12531d6.md -> ffffc29f0024-cadc-4221-9361-18bcb52b968f.md
#   - new xlf correspondence + timestamps for both zh-cn and de-de

$wb = $excel.ActiveWorkbook

$oldUuid1 = "4032afde-8e5a-4dd2-a7f0-f64482d5e798"
$oldUuid2 = "961d70ea-fabc-4d5e-9a78-48ffc12531d6"
$newUuid1 = "c7e37e7b-0198-4224-a154-a931b3cbc493"
$newUuid2 = "ffffc29f0024-cadc-4221-9361-18bcb52b968f"

$newXlfHash = "472148b32b812ac1221de1558083179dd5dd690a"

$newHandbackDate      = "2016-08-22 19:07:26"
$newZhCnHandoffDate   = "2016-08-22 19:07:20"
$newZhCnHandbackDate  = "2016-08-22 19:07:37"
$newDeDeHandbackDate  = "2016-08-22 19:07:44"

$zhCnXlf = "$newUuid1.$newXlfHash.zh-cn.xlf"
$deDeXlf = "$newUuid1.$newXlfHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Sheets.Item("Overview")

$ws1.Range("A2").Value = "$newUuid1.md"
$ws1.Range("A3").Value = "$newUuid2.md"
$ws1.Range("G2").Value = $newHandbackDate
$ws1.Range("G3").Value = $newHandbackDate

# Rebuild the two hyperlinks on B2/B3: the underlying link target is
# unchanged, only the displayed/cell text changes, and re-assigning the
# existing Hyperlink object's properties does not update it in place, so
# drop and recreate both links.
$ws1Link1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid1.md"
$ws1Link2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid2.md"

$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ws1Link1, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newUuid1.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), $ws1Link2, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newUuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Sheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newUuid1.md"
$ws2.Range("I2").Value = "$newUuid1.md"
$ws2.Range("G2").Value = $zhCnXlf
$ws2.Range("J2").Value = $zhCnXlf
$ws2.Range("H2").Value = $newZhCnHandoffDate
$ws2.Range("K2").Value = $newZhCnHandbackDate

$ws2.Range("A3").Value = "$newUuid2.md"
$ws2.Range("I3").Value = "$newUuid2.md"
$ws2.Range("G3").Value = $zhCnXlf
$ws2.Range("J3").Value = $zhCnXlf
$ws2.Range("H3").Value = $newZhCnHandoffDate
$ws2.Range("K3").Value = $newZhCnHandbackDate

$ws2Link1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid1.md"
$ws2Link2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7c85ccfb8dfeb992f7c4d0b6b883aeccb1177453/e2e/$oldUuid1.md"
$ws2Link3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid2.md"
$ws2Link4 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7c85ccfb8dfeb992f7c4d0b6b883aeccb1177453/e2e/$oldUuid2.md"

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2Link1, [System.Type]::Missing, [System.Type]::Missing, "$newUuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $ws2Link2, [System.Type]::Missing, [System.Type]::Missing, "$newUuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2Link3, [System.Type]::Missing, [System.Type]::Missing, "$newUuid2.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $ws2Link4, [System.Type]::Missing, [System.Type]::Missing, "$newUuid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Sheets.Item("de-de")

$ws3.Range("A2").Value = "$newUuid1.md"
$ws3.Range("I2").Value = "$newUuid1.md"
$ws3.Range("G2").Value = $deDeXlf
$ws3.Range("J2").Value = $deDeXlf
$ws3.Range("H2").Value = $newHandbackDate
$ws3.Range("K2").Value = $newDeDeHandbackDate

$ws3.Range("A3").Value = "$newUuid2.md"
$ws3.Range("I3").Value = "$newUuid2.md"
$ws3.Range("G3").Value = $deDeXlf
$ws3.Range("J3").Value = $deDeXlf
$ws3.Range("H3").Value = $newHandbackDate
$ws3.Range("K3").Value = $newDeDeHandbackDate

$ws3Link1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid1.md"
$ws3Link2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bc7dcfe99afeeec49f58ff142b10d7edd08a5716/e2e/$oldUuid1.md"
$ws3Link3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058474ddc5dbe9a8ad595ae082f76724b0b0b00b/e2e/$oldUuid2.md"
$ws3Link4 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bc7dcfe99afeeec49f58ff142b10d7edd08a5716/e2e/$oldUuid2.md"

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3Link1, [System.Type]::Missing, [System.Type]::Missing, "$newUuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $ws3Link2, [System.Type]::Missing, [System.Type]::Missing, "$newUuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ws3Link3, [System.Type]::Missing, [System.Type]::Missing, "$newUuid2.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $ws3Link4, [System.Type]::Missing, [System.Type]::Missing, "$newUuid2.md")

Write-Output "Handback status report regenerated."
